# Updates cryptos list values (price & 1h volume change) per the
# "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value. Cells whose new text parses as a plain
# decimal number (single "." , e.g. "566.70") are written with the cell
# temporarily forced to Text format so Excel does not auto-convert them
# to a Number (which would silently drop the formatting, e.g. "566.70"
# -> 566.7). The format is restored to Normal/General afterwards so no
# visible style change is left behind.
$updates = [ordered]@{
    'D2' = '61.764.69'
    'E2' = '  -1.32%  '
    'D3' = '2.886.99'
    'E3' = '  -2.14%  '
    'E4' = '  +0.02%  '
    'D5' = '566.70'
    'E5' = '  -3.70%  '
    'D6' = '142.88'
    'E6' = '  -2.61%  '
    'E7' = '  +0.03%  '
    'D8' = '0.500'
    'E8' = '  -1.00%  '
    'D9' = '2.884.90'
    'E9' = '  -2.26%  '
    'D10' = '6.95'
    'E10' = '  -0.22%  '
    'E11' = '  -1.88%  '
    'E12' = '  -1.30%  '
    'E13' = '  -0.76%  '
    'D14' = '31.90'
    'E14' = '  -0.90%  '
    'D15' = '0.126'
    'E15' = '  +0.05%  '
    'D16' = '3.366.08'
    'E16' = '  -2.11%  '
    'D17' = '61.719.45'
    'E17' = '  -1.33%  '
    'B18' = 'Polkadot'
    'C18' = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
    'D18' = '6.52'
    'E18' = '  -1.88%  '
    'B19' = 'WrappedEther'
    'C19' = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
    'D19' = '2.880.73'
    'E19' = '  -2.41%  '
    'D20' = '428.22'
    'E20' = '  -1.33%  '
    'D21' = '12.94'
    'E21' = '  -3.86%  '
    'E22' = '  -1.12%  '
    'E23' = '  -1.38%  '
    'D24' = '78.82'
    'E24' = '  -1.63%  '
    'E25' = '  +1.15%  '
    'D26' = '10.02'
    'E26' = '  -9.98%  '
    'E27' = '  +0.07%  '
    'E28' = '  -3.61%  '
    'E29' = '  +9.18%  '
    'D30' = '6.96'
    'E30' = '  -2.76%  '
    'E31' = '  -3.20%  '
    'E32' = '  -6.93%  '
    'E33' = '  +0.02%  '
    'E34' = '  -1.29%  '
    'E35' = '  -2.30%  '
    'E36' = '  -4.44%  '
    'E37' = '  -3.31%  '
    'D38' = '48.84'
    'E39' = '  -6.96%  '
    'E40' = '  -4.80%  '
    'D41' = '0.116'
    'E41' = '  +0.91%  '
    'D42' = '8.12'
    'E42' = '  -2.57%  '
    'D43' = '39.87'
    'E43' = '  +1.65%  '
    'E44' = '  -2.16%  '
    'D45' = '2.688.82'
    'E45' = '  +0.45%  '
    'E46' = '  +0.10%  '
    'D47' = '131.30'
    'E47' = '  -2.73%  '
    'D48' = '344.57'
    'E48' = '  -2.48%  '
    'E49' = '  +0.01%  '
    'E50' = '  -1.37%  '
    'D51' = '21.53'
    'E51' = '  -4.30%  '
}

foreach ($cellRef in $updates.Keys) {
    $value = $updates[$cellRef]
    $range = $ws.Range($cellRef)
    $isPlainNumber = $value.Trim() -match "^[+-]?[0-9]+(\.[0-9]+)?$"
    if ($isPlainNumber) {
        $range.NumberFormat = "@"
        $range.Value = $value
        $range.Style = "Normal"
    } else {
        $range.Value = $value
    }
}
